$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# 1. Re-create the "smarthosting" bookmark so Word re-issues it a fresh
#    (lower) internal id when the document is saved.
$bm = $d.Bookmarks("smarthosting")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("smarthosting", $bmRange)

# 2. Remove the CryptoBridge and HitBTC HYPERLINK field codes entirely.
$d.Fields.Item(1).Delete()
$d.Fields.Item(1).Delete()

# 3. Trim " such as" (plus the trailing nbsp) from the sentence, leaving
#    "... obtained from exchanges" directly followed by the final period.
$r1 = $d.Content
$r1.Find.Execute(" such as$nbsp", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 4. Remove the now-orphaned ", " connector run that used to sit between
#    the two removed hyperlinks.
$r2 = $d.Content
$r2.Find.Execute(",$nbsp", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
